$wb = $excel.ActiveWorkbook

# Data describing each sheet: numeric tail (already present in col A before
# the shift, kept as numbers in col B) and the new file-name strings that
# get appended below them (as strings in col B).
$sheetsData = @(
    @{
        Index = 1
        Numbers = @(1, 2, 3)
        Files = @("file1.pdf", "file2.pdf", "file3.pdf", "file4.pdf")
    },
    @{
        Index = 2
        Numbers = @(3, 4, 5)
        Files = @("file6.pdf", "file4.pdf", "file5.pdf")
    },
    @{
        Index = 3
        Numbers = @(6, 7, 8)
        Files = @(
            "file5 copia 2.pdf",
            "file5 copia.pdf",
            ".DS_Store",
            "file4 copia.pdf",
            "file6.pdf",
            "file6 copia 2.pdf",
            "file4.pdf",
            "file5.pdf",
            "file4 copia 2.pdf",
            "file6 copia.pdf"
        )
    }
)

foreach ($sd in $sheetsData) {
    $ws = $wb.Worksheets.Item($sd.Index)

    # Insert a fresh column before A: old col A (header + numbers) shifts to B.
    $ws.Columns.Item(1).Insert()

    $headerCell = $ws.Range("B1")

    $row = 2
    $idx = 0

    # Rows that already held the numeric values (now in column B); add the
    # running index 0,1,2,... in the new column A alongside them.
    foreach ($n in $sd.Numbers) {
        $aCell = $ws.Range("A" + $row)
        $aCell.Value = $idx
        $headerCell.Copy()
        $aCell.PasteSpecial(-4122)
        $row = $row + 1
        $idx = $idx + 1
    }

    # New rows: index in column A, file name string in column B.
    foreach ($f in $sd.Files) {
        $aCell = $ws.Range("A" + $row)
        $bCell = $ws.Range("B" + $row)
        $aCell.Value = $idx
        $headerCell.Copy()
        $aCell.PasteSpecial(-4122)
        $bCell.Value = $f
        $row = $row + 1
        $idx = $idx + 1
    }
}

$excel.CutCopyMode = $false
